# Add 7 new rows (160-166) of data to Sheet1, continuing the existing
# time series that currently ends at row 159 (date serial 45589).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then a map of column letter -> value for the
# columns that are populated in this data set (A, B, C, I, K, N, O, Q, U, Z).
# Note: values that were written in scientific notation in the source data
# (column Q) are expanded to plain decimal literals below since the script
# parser does not accept exponent notation (e.g. "2.1072E-06").
$newRows = @(
    @{ Row = 160; A = 45590; B = 590.3322451974001;  C = 169.223682537;   I = 295.24406719;   K = 78.924145584798;  N = 31.81545294816;  O = 0.9794309909999999; Q = 0.0000021072; U = 218.1189759152305; Z = 592.8770018009119 },
    @{ Row = 161; A = 45591; B = 593.8232583528001;  C = 172.1281822385;  I = 305.950916928;  K = 80.18615529279602; N = 32.09484221664; O = 0.993220998;         Q = 0.0000021792; U = 219.3982660965515; Z = 589.9424283319099 },
    @{ Row = 162; A = 45592; B = 602.045101926;      C = 173.88169853;    I = 316.245964753;  K = 79.60368927372001; N = 32.57213221696; O = 1.003095324;         Q = 0.0000022008; U = 226.434362093817;  Z = 684.7338094337999 },
    @{ Row = 163; A = 45593; B = 619.2201290238002;  C = 178.019691898;   I = 319.218100299;  K = 79.70076694356601; N = 32.54884977792; O = 1.023354717;         Q = 0.0000022032; U = 232.3190969278936; Z = 632.6567754918281 },
    @{ Row = 164; A = 45594; B = 643.7740514076;     C = 182.96476038;    I = 321.241301002;  K = 80.720082476949;   N = 34.32995636448; O = 1.032207561;         Q = 0.0000023208; U = 244.7282116867073; Z = 648.1215118681559 },
    @{ Row = 165; A = 45595; B = 640.3073778972001;  C = 184.3785285565;  I = 313.094784897;  K = 82.127708689716;   N = 33.78281904704; O = 1.013991132;         Q = 0.0000023136; U = 240.6344831064801; Z = 666.3810801197241 },
    @{ Row = 166; A = 45596; B = 622.1391162678;     C = 174.6312244735;  I = 302.029846539;  K = 81.739398010332;   N = 31.62919343584; O = 0.981644202;         Q = 0.0000021864; U = 226.1785040575528; Z = 652.546662337286 }
)

$colIndex = @{
    A = 1; B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8; I = 9; J = 10;
    K = 11; L = 12; M = 13; N = 14; O = 15; P = 16; Q = 17; R = 18; S = 19;
    T = 20; U = 21; V = 22; W = 23; X = 24; Y = 25; Z = 26
}

# The last existing data row (159) has cell A159 styled with the date
# number format used throughout column A. Copy that formatting onto each
# new A-column cell so the new rows match the existing look.
$ws.Cells.Item(159, 1).Copy() | Out-Null

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    foreach ($key in $entry.Keys) {
        if ($key -eq "Row") { continue }
        $colNum = $colIndex[$key]
        $ws.Cells.Item($r, $colNum).Value = $entry[$key]
    }
}

$excel.CutCopyMode = 0
